$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M1").Value = 'imagen'
$ws.Range("M2").Value = 'https://fernandobelaundeterry.com.pe/img/galeria/e3565820e160444c5a5bc5c61b58148e83c54bb7.jpg'
$ws.Range("M3").Value = 'https://pbs.twimg.com/media/GN0pfFIW0AACKeg.png'
$ws.Range("M4").Value = 'https://www.servindi.org/sites/default/files/editor/imagenes/gein.png'
$ws.Range("M5").Value = 'https://wayka.pe/wp-content/uploads/2017/01/uchu.jpg.webp'
$ws.Range("M6").Value = 'https://cdn01.pucp.education/idehpucp/images/2024/03/imagen_2022-07-19_063244840.png'
$ws.Range("M7").Value = 'https://media.cnn.com/api/v1/images/stellar/prod/cnne-638827-peru-politics-garcia.jpg?c=original'
$ws.Range("M8").Value = 'https://lineadetiempo.iep.org.pe/public/upload/file_569fc11cddd8d.jpg'
$ws.Range("M9").Value = 'https://pbs.twimg.com/media/Evhz268WYAMTBsI.jpg'
$ws.Range("M10").Value = 'https://scontent.flim28-1.fna.fbcdn.net/v/t1.6435-9/34963073_10156190722935661_1021843668283162624_n.png?stp=dst-png_p526x296&_nc_cat=106&ccb=1-7&_nc_sid=127cfc&_nc_ohc=gMMEfvlABJIQ7kNvwHw4t2x&_nc_oc=AdmeQiKQt8Trz1mHjtlIqhAm7FW8j-ttjqSUyzsOtaaGuPmbSWbjfVfvZ1dQwUL1LFE&_nc_zt=23&_nc_ht=scontent.flim28-1.fna&_nc_gid=3-njooDcKe5ZrizytwiXGQ&oh=00_AfksBDAzJSmZcf9hSuYfB04WbPDI_wwRxLKzF6doZCb4uA&oe=695BBF63'
$ws.Range("M11").Value = 'https://cloudfront-us-east-1.images.arcpublishing.com/infobae/MNMIX7O6JNGI3IGBS7ZZE54AAQ.png'
$ws.Range("M12").Value = 'https://s.france24.com/media/display/b4652ef0-e9fe-11eb-9ff8-005056a97e36/w:1280/p:16x9/294af33ffc0dcaa192c1a3a98139bb16db45c844.jpg'
$ws.Range("M13").Value = 'https://images.ctfassets.net/owmn3xt273ds/1Ta2J1EbrBcDECGkGcgdQf/4d4d201c65b0ceb4f8a73fc4abce10e1/NS231120_Caballero.jpg?fm=webp'
$ws.Range("M14").Value = 'https://egecperu.pe/wp-content/uploads/2025/02/El-Conflicto-del-Cenepa-La-ultima-guerra-de-America-Latina-y-su-legado-30-anos-despues.webp'
$ws.Range("M15").Value = 'https://imagenes.elpais.com/resizer/v2/IFSHOUIXSUU7AG2L7F4Y6JGX5M.jpg?auth=d9163fb9b90528567b5431808e2a54298e3cc1bc2e283c3a8f655e9e420319ca&width=1960&height=1103&smart=true'
$ws.Range("M16").Value = 'https://lum.cultura.pe/cdi/sites/default/files/articulo/fotos/3_5.jpg'
$ws.Range("M20").Value = 'https://cloudfront-us-east-1.images.arcpublishing.com/infobae/UL6SX2MWK5DCJFDTA63U74BCCI.jpg'
$ws.Range("M18").Value = 'https://imagenes.elpais.com/resizer/v2/B2SC7E7AYJXRZ76B3V3KFHLV4Y.jpg?auth=849475f022e746d9bde24b95cf684adcb9a3e97b3ee2d47b5ad1c406688731be&width=414'
$ws.Range("M19").Value = 'https://ceeep.mil.pe/wp-content/uploads/2023/05/word-image-2.png'
$ws.Range("M17").Value = 'https://www.infobae.com/new-resizer/x6tcVNbY5McXHQE1WcwLnQvwKSY=/arc-anglerfish-arc2-prod-infobae/public/QWXC4RD3EZE2DC2YPG2LWVNJ6Q.jpg'
$ws.Range("M21").Value = 'https://www.proetica.org.pe/wp-content/uploads/2018/01/LAVA-JATO2.jpg'
$ws.Range("M22").Value = 'https://larepublica.cronosmedia.glr.pe/original/2022/05/14/627f70f8580a3a2b71793e9f.jpg'
$ws.Range("M23").Value = 'https://cde.canaln.pe/actualidad-gobierno-alan-garcia-percibido-como-mas-corrupto-ultimos-anos-n293761-696x418-405589.jpg'
$ws.Range("M24").Value = 'https://www.infobae.com/new-resizer/d6B9iKFMlztidbjW1KNsWUXRJSI=/arc-anglerfish-arc2-prod-infobae/public/ESMWDG2KPNAOXOEXPVPA47O3HA.png'
$ws.Range("M25").Value = 'https://portal.andina.pe/EDPfotografia2/Thumbnail/2009/03/10/000089188W.jpg'
$ws.Range("M26").Value = 'https://imgs.mongabay.com/wp-content/uploads/sites/25/2016/09/03171652/bagua_conflicto_pueblo-indigena_mongabay-1.jpg'
$ws.Range("M27").Value = 'https://lineadetiempo.iep.org.pe/public/upload/file_569955f6af2a9.jpg'
$ws.Range("M28").Value = 'https://upload.wikimedia.org/wikipedia/commons/3/30/Ollanta_Humala_asume_Presidencia_del_Per%C3%BA.jpg'
$ws.Range("M29").Value = 'https://lineadetiempo.iep.org.pe/public/upload/file_5699ae8de29a7.jpg'
$ws.Range("M30").Value = 'https://scontent.flim38-1.fna.fbcdn.net/v/t1.6435-9/167613447_5263535013687951_4723074474013864240_n.jpg?_nc_cat=108&ccb=1-7&_nc_sid=127cfc&_nc_ohc=o8tz-1k1JycQ7kNvwHkUk4w&_nc_oc=Adk6icTyWnS7dsPJHdEAFMSHEiptTP0NaeuVmLAnCZodlQs_71Izs5NWgLLZiLxiLxk&_nc_zt=23&_nc_ht=scontent.flim38-1.fna&_nc_gid=LxrPtW7f46WvpCpLpaCVtg&oh=00_AfnnFBboSYUyQG-ljdYywE2ox1Gp_D8RD1EmKfQD2wSeBg&oe=695C0000'
$ws.Range("M31").Value = 'https://cdn01.pucp.education/idehpucp/images/2024/03/imagen_2022-09-26_232350211-1024x683.png'
$ws.Range("M32").Value = 'https://upload.wikimedia.org/wikipedia/commons/4/43/MapaElectoralRegionalPeru2014.png'
$ws.Range("M33").Value = 'https://elperuano.pe/suplementosflipping/lonuestro/251/web/img/foto06.jpg'
$ws.Range("M34").Value = 'https://cloudfront-eu-central-1.images.arcpublishing.com/prisaradio/GXMOG2MOFNLXRLNZUDU24WINMY.jpg'
$ws.Range("M35").Value = 'https://museos.cultura.pe/sites/default/files/noticias/imagenes/2021-12/Museo%20Tu%CC%81cume.png'
$ws.Range("M36").Value = 'https://ojo-publico.com/sites/default/files/32701602463_d00f89db08_h.jpg'
$ws.Range("M37").Value = 'https://www.infobae.com/resizer/v2/V7TPPRJTPBENLBRP52G7DVCY2Y.png?auth=ed6944ef06a5f942357505fc80ca93aded1e766dda03d0a2dc9a7fd8dc218e06&smart=true&width=1200&height=1200&quality=85'
$ws.Range("M38").Value = 'https://img.lpderecho.pe/wp-content/uploads/2021/07/Licencia-sindical-Pedro-Castillo-LP.jpg'
$ws.Range("M39").Value = 'https://ucsp.edu.pe/wp-content/uploads/2023/08/b1161.jpeg'
$ws.Range("M40").Value = 'https://scontent.flim38-1.fna.fbcdn.net/v/t39.30808-6/564716127_1212449430920169_8072068678252789608_n.jpg?_nc_cat=105&ccb=1-7&_nc_sid=833d8c&_nc_ohc=EdDd3VlrdeYQ7kNvwGJISI9&_nc_oc=AdnUiTEhUP3yz_ZfLs8DBk6_NP_bN_aH4yF7vuTeKrLfx3VJAwH5gtETNr0QJrWLHl8&_nc_zt=23&_nc_ht=scontent.flim38-1.fna&_nc_gid=S34O2x7lt3kRWQgKaxrrTw&oh=00_Afkr5Sbi5upttlaBB4bMTkFG2m5lSbKJj3waESqPSETzQA&oe=693A4AC9'
$ws.Range("M41").Value = 'https://peru.info/archivos/publicacion/26-imagen-2035341682017.jpg'
$ws.Range("M42").Value = 'https://cooperaccion.org.pe/wp-content/uploads/2019/11/AYACUCHO-PARO-MINER%C3%8DA-La-mula.jpg'
$ws.Range("M43").Value = 'https://i.ytimg.com/vi/Fdp2kEMhG7g/maxresdefault.jpg'
$ws.Range("M44").Value = 'https://www.infobae.com/new-resizer/YxvTlJzVnFWweLyB8sUyhdL2n5c=/arc-anglerfish-arc2-prod-infobae/public/DCUTTZYEXXKH3OKVMGZKOBTNTU.jpg'
$ws.Range("M45").Value = 'https://upload.wikimedia.org/wikipedia/commons/8/82/Protestas_en_Lima_por_golpe_de_estado_parlamentario_27.jpg'
$ws.Range("M46").Value = 'https://www.unicef.org/peru/sites/unicef.org.peru/files/styles/press_release_feature/public/Vacunacion%20nino%20Peru%20-%20UNICEF%20UN0724497.jpg.webp?itok=abBnc3U0'
